$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4:D4").Copy() | Out-Null
$ws.Range("A5:D5").PasteSpecial() | Out-Null

$ws.Range("A5").Value = "Carlos"
$ws.Range("B5").Value = "carlos"
$ws.Range("C5").Value = "2021002252@ifam.edu.br"
$ws.Range("D5").Value = "admin123"
